$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (new report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 32   Number  11"
$ws.Range("C9").Value = "Report Covering the Week  3/10/2025  Through  3/16/2025"

# --- Row 15 (Murder) : only the 2-Year / 15-Year %chg columns move ---
$ws.Range("M15").Value = -83.333333333333
$ws.Range("N15").Value = -92.857142857142

# --- Row 16 (Rape) ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = -60
$ws.Range("I16").Value = 17
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -37.037037037037
$ws.Range("L16").Value = -19.047619047619
$ws.Range("M16").Value = -51.428571428571
$ws.Range("N16").Value = -89.634146341463

# --- Row 17 (Robbery) ---
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 32
$ws.Range("J17").Value = 37
$ws.Range("K17").Value = -13.513513513513
$ws.Range("L17").Value = -25.581395348837
$ws.Range("M17").Value = 33.333333333333
$ws.Range("N17").Value = -74.603174603174

# --- Row 18 (Fel. Assault) ---
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 15
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -44.444444444444
$ws.Range("M18").Value = -28.571428571428
$ws.Range("N18").Value = -94.252873563218

# --- Row 19 (Burglary) ---
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 19
$ws.Range("G19").Value = 24
$ws.Range("H19").Value = -20.833333333333
$ws.Range("I19").Value = 57
$ws.Range("J19").Value = 78
$ws.Range("K19").Value = -26.923076923076
$ws.Range("L19").Value = 1.785714285714
$ws.Range("M19").Value = 119.230769230769
$ws.Range("N19").Value = -28.75

# --- Row 20 (Gr. Larceny) ---
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 4
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 11
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -35.294117647058
$ws.Range("L20").Value = -47.619047619047
$ws.Range("M20").Value = 22.222222222222
$ws.Range("N20").Value = -80

# --- Row 21 (G.L.A. - bold total-ish row) ---
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -38.888888888888
$ws.Range("F21").Value = 46
$ws.Range("G21").Value = 62
$ws.Range("H21").Value = -25.806451612903
$ws.Range("I21").Value = 128
$ws.Range("J21").Value = 177
$ws.Range("K21").Value = -27.683615819209
$ws.Range("L21").Value = -20.496894409937
$ws.Range("M21").Value = 11.304347826087
$ws.Range("N21").Value = -79.552715654952

# --- Row 22 (TOTAL) : C/D/E/F flip from "N/A" placeholders to real numbers ---
$ws.Range("C22").Value = 1
$ws.Range("C22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;`"-`"#,##0.0"
$ws.Range("F22").Value = 1
$ws.Range("F22").NumberFormat = "#,##0"
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = -20
$ws.Range("L22").Value = -60
$ws.Range("M22").Value = 33.333333333333

# --- Row 23 (Transit) ---
$ws.Range("L23").Value = -75

# --- Row 24 (Housing) ---
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = -57.894736842105
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = 7.936507936507
$ws.Range("I24").Value = 146
$ws.Range("J24").Value = 166
$ws.Range("K24").Value = -12.048192771084
$ws.Range("L24").Value = -7.006369426751
$ws.Range("M24").Value = 160.714285714286

# --- Row 25 (Petit Larceny) ---
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -50
$ws.Range("G25").Value = 12
$ws.Range("H25").Value = 66.666666666666
$ws.Range("I25").Value = 28
$ws.Range("J25").Value = 30
$ws.Range("K25").Value = -6.666666666666
$ws.Range("L25").Value = -30

# --- Row 26 (Retail Theft) ---
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 125
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 13.793103448275
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 62
$ws.Range("K26").Value = 11.290322580645
$ws.Range("L26").Value = 18.965517241379
$ws.Range("M26").Value = -22.471910112359

# --- Row 28 (UCR Rape*) : C28 flips from "N/A" placeholder to a real number ---
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 2
$ws.Range("K28").Value = -66.666666666666
$ws.Range("L28").Value = -50

# --- Rows 29 / 30 (Shooting Vic. / Shooting Inc.) : 2-Year %chg refresh ---
$ws.Range("N29").Value = -95.652173913043
$ws.Range("N30").Value = -95.652173913043
